$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("B3").Value = 2.84720710080662
$ws.Range("F3").Value = 2.72290869198426
$ws.Range("G3").Value = 1.9957455552218799
$ws.Range("J3").Value = 5.2154636703621504
$ws.Range("K3").Value = 5.0521390309110803

# Row 4
$ws.Range("B4").Value = 2.84552241201217
$ws.Range("C4").Value = 1.81907116384734
$ws.Range("F4").Value = 2.9758989987567599
$ws.Range("G4").Value = 2.0597339977487601
$ws.Range("J4").Value = 17.0824840671484
$ws.Range("K4").Value = 13.806472688861801

# Row 5
$ws.Range("B5").Value = 2.1875880758172599
$ws.Range("C5").Value = 1.69557812916897
$ws.Range("F5").Value = 2.56529643278509
$ws.Range("J5").Value = 29.299322874224401
$ws.Range("K5").Value = 18.243605474790002

# Row 6
$ws.Range("B6").Value = 2.18918349455274
$ws.Range("C6").Value = 1.6952580503544199
$ws.Range("G6").Value = 1.88709355377798
$ws.Range("J6").Value = 73.429297029762395
$ws.Range("K6").Value = 39.420138453430702

# Row 12
$ws.Range("B12").Value = 1.7072244343728999
$ws.Range("C12").Value = 1.9148926318613699

# Row 22
$ws.Range("B22").Value = 2.04
$ws.Range("C22").Value = 1.42828568570857
$ws.Range("F22").Value = 3.11944768859051
$ws.Range("G22").Value = 2.0770919684817999
$ws.Range("J22").Value = 3.6232989534435398
$ws.Range("K22").Value = 2.3357173877071098

# Row 23
$ws.Range("B23").Value = 2.2200000000000002
$ws.Range("C23").Value = 1.6572660052764301
$ws.Range("F23").Value = 2.2917008817116602
$ws.Range("G23").Value = 1.57659259054217
$ws.Range("J23").Value = 3.3612493075810099
$ws.Range("K23").Value = 2.4016888955013802

# Row 24
$ws.Range("B24").Value = 2.2216144420414801
$ws.Range("C24").Value = 1.66388687383708
$ws.Range("F24").Value = 2.4095408813260901
$ws.Range("G24").Value = 1.6677210442796699
$ws.Range("J24").Value = 3.4175288841464901
$ws.Range("K24").Value = 2.8849892408554001

# Row 25
$ws.Range("B25").Value = 2.5675539528999001
$ws.Range("C25").Value = 1.70267854864918
$ws.Range("F25").Value = 2.2350394903862298
$ws.Range("G25").Value = 1.72692469676862
$ws.Range("J25").Value = 2.74586047407964
$ws.Range("K25").Value = 2.3702314644120199

# Move the selection to A1 (the saved sheet no longer has a B3 selection).
$ws.Range("A1").Select()
